$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42:151 down to 43:152
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with its data
$ws.Cells.Item(42, 1).Value = 1
$ws.Cells.Item(42, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(42, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(42, 4).Value = 44498
$ws.Cells.Item(42, 5).Value = 15
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100108
$ws.Cells.Item(42, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(42, 9).Value = 100108006
$ws.Cells.Item(42, 10).Value = "Plátano"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Pintón"
$ws.Cells.Item(42, 13).Value = 120
$ws.Cells.Item(42, 14).Value = 19000
$ws.Cells.Item(42, 15).Value = 20000
$ws.Cells.Item(42, 16).Value = 19500
$ws.Cells.Item(42, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(42, 18).Value = "Ecuador"
$ws.Cells.Item(42, 19).Value = 975
$ws.Cells.Item(42, 20).Value = 20
